# Auto-generated PowerShell COM-interop script
# Applies the scheduled-runner market-data update to Halicarnassus_Profits
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 803.25
$ws.Range("I6").Value = 57
$ws.Range("K6").Value = 171
$ws.Range("M6").Value = -59

$ws.Range("H9").Value = 172.83333
$ws.Range("I9").Value = 184.25
$ws.Range("J9").Value = 150
$ws.Range("K9").Value = 184.25
$ws.Range("L9").Value = 150
$ws.Range("M9").Value = -15.25
$ws.Range("N9").Value = -488

$ws.Range("H62").Value = 5204.1577
$ws.Range("I62").Value = 2180
$ws.Range("K62").Value = 2180
$ws.Range("M62").Value = -1556

$ws.Range("H65").Value = 5204.1577
$ws.Range("I65").Value = 2180
$ws.Range("K65").Value = 10900
$ws.Range("M65").Value = -7780

$ws.Range("H94").Value = 3925
$ws.Range("I94").Value = 4138.8887
$ws.Range("K94").Value = 4138.8887
$ws.Range("M94").Value = -3687.8887

$ws.Range("H98").Value = 1080.1428
$ws.Range("I98").Value = 1038.75
$ws.Range("K98").Value = 1038.75
$ws.Range("M98").Value = 459.25

$ws.Range("H122").Value = 1080.1428
$ws.Range("I122").Value = 1038.75
$ws.Range("K122").Value = 3116.25
$ws.Range("M122").Value = -666.25

$ws.Range("H135").Value = 1173.1904
$ws.Range("I135").Value = 948.55554
$ws.Range("J135").Value = 2521
$ws.Range("K135").Value = 8536.99986
$ws.Range("L135").Value = 22689
$ws.Range("M135").Value = -6001.99986
$ws.Range("N135").Value = -27759


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5766.7144
$ws.Range("I2").Value = 4690.8
$ws.Range("K2").Value = 4690.8
$ws.Range("M2").Value = -4577.8

$ws.Range("H61").Value = 1981.8182
$ws.Range("I61").Value = 2062.5
$ws.Range("J61").Value = 1766.6666
$ws.Range("K61").Value = 2062.5
$ws.Range("L61").Value = 1766.6666
$ws.Range("M61").Value = -1850.5
$ws.Range("N61").Value = -2190.6666

$ws.Range("H116").Value = 5766.7144
$ws.Range("I116").Value = 4690.8
$ws.Range("K116").Value = 4690.8
$ws.Range("M116").Value = -2396.8

$ws.Range("H132").Value = 2949.8333
$ws.Range("I132").Value = 2946.8823
$ws.Range("K132").Value = 8840.6469
$ws.Range("M132").Value = -6310.6469

$ws.Range("H136").Value = 1981.8182
$ws.Range("I136").Value = 2062.5
$ws.Range("J136").Value = 1766.6666
$ws.Range("K136").Value = 6187.5
$ws.Range("L136").Value = 5299.9998
$ws.Range("M136").Value = -3637.5
$ws.Range("N136").Value = -10399.9998


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5766.7144
$ws.Range("I3").Value = 4690.8
$ws.Range("K3").Value = 4690.8
$ws.Range("M3").Value = -4576.8

$ws.Range("H86").Value = 3738.8572
$ws.Range("I86").Value = 1603.3334
$ws.Range("J86").Value = 7582.8
$ws.Range("K86").Value = 1603.3334
$ws.Range("L86").Value = 7582.8
$ws.Range("M86").Value = -480.3334
$ws.Range("N86").Value = -9828.799999999999

$ws.Range("H88").Value = 1294681.9
$ws.Range("J88").Value = 1294681.9
$ws.Range("L88").Value = 1294681.9
$ws.Range("N88").Value = -1295493.9

$ws.Range("H89").Value = 3738.8572
$ws.Range("I89").Value = 1603.3334
$ws.Range("J89").Value = 7582.8
$ws.Range("K89").Value = 8016.666999999999
$ws.Range("L89").Value = 37914
$ws.Range("M89").Value = -2400.666999999999
$ws.Range("N89").Value = -49146

$ws.Range("H91").Value = 1294681.9
$ws.Range("J91").Value = 1294681.9
$ws.Range("L91").Value = 1294681.9
$ws.Range("N91").Value = -1297489.9

$ws.Range("H130").Value = 62955.6
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 62955.6
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 62955.6
$ws.Range("M130").ClearContents()
$ws.Range("N130").Value = -72995.60000000001


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1164.1333
$ws.Range("I107").Value = 948.4
$ws.Range("J107").Value = 1595.6
$ws.Range("K107").Value = 948.4
$ws.Range("L107").Value = 1595.6
$ws.Range("M107").Value = 971.6
$ws.Range("N107").Value = -5435.6


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 53958.332
$ws.Range("I11").Value = 240000
$ws.Range("J11").Value = 5000
$ws.Range("K11").Value = 720000
$ws.Range("L11").Value = 15000
$ws.Range("M11").Value = -719860
$ws.Range("N11").Value = -15280

$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()

$ws.Range("H122").Value = 812.2222
$ws.Range("I122").Value = 851.7143
$ws.Range("J122").Value = 674
$ws.Range("K122").Value = 7665.428699999999
$ws.Range("L122").Value = 6066
$ws.Range("M122").Value = -5215.428699999999
$ws.Range("N122").Value = -10966

$ws.Range("H138").Value = 4566
$ws.Range("I138").Value = 1515
$ws.Range("J138").Value = 6600
$ws.Range("K138").Value = 4545
$ws.Range("L138").Value = 19800
$ws.Range("M138").Value = 595
$ws.Range("N138").Value = -30080


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 3691.2
$ws.Range("I31").Value = 2114
$ws.Range("K31").Value = 2114
$ws.Range("M31").Value = -1822

$ws.Range("H37").Value = 3691.2
$ws.Range("I37").Value = 2114
$ws.Range("K37").Value = 2114
$ws.Range("M37").Value = -1837

$ws.Range("H132").Value = 105343.6
$ws.Range("I132").Value = 171256.33
$ws.Range("K132").Value = 513768.99
$ws.Range("M132").Value = -511238.99


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6228.5713
$ws.Range("I61").Value = 4300
$ws.Range("J61").Value = 7000
$ws.Range("K61").Value = 4300
$ws.Range("L61").Value = 7000
$ws.Range("M61").Value = -4098
$ws.Range("N61").Value = -7404

$ws.Range("H106").Value = 12499.25
$ws.Range("J106").Value = 12499.25
$ws.Range("L106").Value = 12499.25
$ws.Range("N106").Value = -15023.25

$ws.Range("H113").Value = 6228.5713
$ws.Range("I113").Value = 4300
$ws.Range("J113").Value = 7000
$ws.Range("K113").Value = 4300
$ws.Range("L113").Value = 7000
$ws.Range("M113").Value = -2130
$ws.Range("N113").Value = -11340


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4333
$ws.Range("I132").Value = 3999.5
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 11998.5
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -9468.5
$ws.Range("N132").Value = -20060

